$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-52 down to 3-53 (bottom-up, cell-by-cell to preserve
# numeric types/styles and avoid Excel auto-cloning row-insert formatting).
for ($r = 52; $r -ge 2; $r--) {
    for ($c = 1; $c -le 5; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r + 1, $c)
        $dstCell.Value2 = $srcCell.Value2
    }
}

# Row 53 is newly created beyond the original range; copy the date style from A52.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = $False

# Populate the new row 2 with new forecast data.
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 5.896808312953783
$ws.Range("D2").Value = 2008

# Update recomputed y_0_forecast (C) and y_1_forecast (E) values for shifted rows.
$ws.Range("C4").Value = 7.441962824572235
$ws.Range("C6").Value = 6.277541464866987
$ws.Range("C7").Value = 5.907218141265402
$ws.Range("E7").Value = 5.917486466529609
$ws.Range("C8").Value = 6.535114773304773
$ws.Range("E8").Value = 6.325696408067327
$ws.Range("C9").Value = 6.511263427347003
$ws.Range("E9").Value = 6.422943767670297
$ws.Range("C10").Value = 5.12051970717502
$ws.Range("E10").Value = 4.950888348161886
$ws.Range("C11").Value = 4.325828829470257
$ws.Range("E11").Value = 4.493586323244281
$ws.Range("C12").Value = 3.65682115264816
$ws.Range("E12").Value = 3.982564147794321
$ws.Range("C13").Value = 2.964652118442834
$ws.Range("E13").Value = 3.83627393798931
$ws.Range("C14").Value = 2.943878639034381
$ws.Range("E14").Value = 4.334309403335435
$ws.Range("C15").Value = 3.171852776411788
$ws.Range("E15").Value = 4.030605385534614
$ws.Range("C16").Value = 1.172679597477866
$ws.Range("E16").Value = 2.644356903452572
$ws.Range("C17").Value = 1.773712379859993
$ws.Range("E17").Value = 3.484530515673856
$ws.Range("C18").Value = 2.961845079861303
$ws.Range("E18").Value = 3.383932287548697
$ws.Range("C19").Value = 2.533350906619081
$ws.Range("E19").Value = 3.524103740130435
$ws.Range("C20").Value = 2.508469427909898
$ws.Range("E20").Value = 3.355044026998955
$ws.Range("C21").Value = 2.661040979345697
$ws.Range("E21").Value = 3.567108445582057
$ws.Range("C22").Value = 3.523703831572056
$ws.Range("E22").Value = 3.74984170812418
$ws.Range("C23").Value = 2.788213251109939
$ws.Range("E23").Value = 3.595252063027843
$ws.Range("C24").Value = 3.150198973767537
$ws.Range("E24").Value = 3.699072253610103
$ws.Range("C25").Value = 1.334931695392405
$ws.Range("E25").Value = 2.689909849380556
$ws.Range("C26").Value = 1.178605266817589
$ws.Range("E26").Value = 2.186196327763934
$ws.Range("C27").Value = 1.80511617406458
$ws.Range("E27").Value = 3.037278170871094
$ws.Range("C28").Value = 0.4641929091049102
$ws.Range("E28").Value = 2.550259844884462
$ws.Range("C29").Value = 2.9927258084951
$ws.Range("E29").Value = 2.807906319450781
$ws.Range("C30").Value = 3.047037961814492
$ws.Range("E30").Value = 2.880436144359444
$ws.Range("C31").Value = 3.113514644866355
$ws.Range("E31").Value = 2.798268189979214
$ws.Range("C32").Value = 2.585454129751663
$ws.Range("E32").Value = 2.671828487424377
$ws.Range("C33").Value = -0.4891791466461126
$ws.Range("E33").Value = 1.156476476708135
$ws.Range("C34").Value = -0.2228847697281378
$ws.Range("E34").Value = 1.982741503124119
$ws.Range("C35").Value = 2.240953541724267
$ws.Range("E35").Value = 2.598498189609066
$ws.Range("C36").Value = -0.4532848472497908
$ws.Range("E36").Value = 2.066462658785673
$ws.Range("C37").Value = -0.7941560676977599
$ws.Range("E37").Value = 2.320050994894562
$ws.Range("C38").Value = -1.165854108406617
$ws.Range("E38").Value = 2.782217648649521
$ws.Range("C39").Value = 3.392010093835562
$ws.Range("E39").Value = 2.595837839692172
$ws.Range("C40").Value = 0.9477102747197819
$ws.Range("E40").Value = 1.83067479293082
$ws.Range("C41").Value = 2.34069710769782
$ws.Range("E41").Value = 2.116186509693896
$ws.Range("C42").Value = 2.501311189006916
$ws.Range("E42").Value = 2.985901060752827
$ws.Range("C43").Value = 2.735256324140778
$ws.Range("E43").Value = 2.169811364059249
$ws.Range("C44").Value = 2.222852754198135
$ws.Range("E44").Value = 1.934107558751452
$ws.Range("C45").Value = 0.9259311313598806
$ws.Range("E45").Value = 1.280651803859989
$ws.Range("C46").Value = 0.6753076481029074
$ws.Range("E46").Value = 0.7957830962485257
$ws.Range("C47").Value = 1.324987171138314
$ws.Range("E47").Value = 1.882659757536698
$ws.Range("C48").Value = -0.007094633234694392
$ws.Range("E48").Value = 1.444584248586422
$ws.Range("C49").Value = 1.935025917091848
$ws.Range("E49").Value = 2.004076350201744
$ws.Range("C50").Value = 2.039329803030121
$ws.Range("E50").Value = 2.510359031091491
$ws.Range("C51").Value = 1.751699155751707
$ws.Range("E51").Value = 1.973546706924467
$ws.Range("C52").Value = 3.078872076370009
$ws.Range("E52").Value = 2.279508996785351

# Populate the new final row 53 with new forecast data.
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = 2.43119486791763
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 2.421949074001883
